$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.552.87"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.93"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.61"
$ws.Range("E5").Value = "  +5.25%  "
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.24"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.337"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.184.28"
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.44"
$ws.Range("E13").Value = "  +8.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.909.99"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.84"
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.518.28"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.92"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0810"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.78"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.90"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.86"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").Value = "  +24.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.56"
$ws.Range("E27").Value = "  +8.52%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +33.02%  "
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.14"
$ws.Range("E34").Value = "  +5.22%  "
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.11"
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("B39").Value = "MultiversX"
$ws.Range("C39").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "53.46"
$ws.Range("E39").Value = "  +56.08%  "
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "91.37"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.354.70"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.46"
$ws.Range("E43").Value = "  +5.93%  "
$ws.Range("E44").Value = "  +11.84%  "
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.72"
$ws.Range("E46").Value = "  +7.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.65"
$ws.Range("E49").Value = "  +5.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.091.55"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0689"
$ws.Range("E51").Value = "  +2.60%  "
